$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3796.0715
$ws.Range("I40").Value = 2295.4285
$ws.Range("K40").Value = 2295.4285
$ws.Range("M40").Value = -2120.4285
$ws.Range("H86").Value = 4649.5
$ws.Range("J86").Value = 5567.3335
$ws.Range("L86").Value = 5567.3335
$ws.Range("N86").Value = -7813.3335
$ws.Range("H89").Value = 4649.5
$ws.Range("J89").Value = 5567.3335
$ws.Range("L89").Value = 27836.6675
$ws.Range("N89").Value = -39068.6675
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").ClearContents()
$ws.Range("H125").Value = 2377.1667
$ws.Range("I125").Value = 972.6
$ws.Range("J125").Value = 9400
$ws.Range("K125").Value = 8753.4
$ws.Range("L125").Value = 84600
$ws.Range("M125").Value = -6293.4
$ws.Range("N125").Value = -89520
$ws.Range("H132").Value = 13350.1455
$ws.Range("I132").Value = 1897.5122
$ws.Range("K132").Value = 5692.536599999999
$ws.Range("M132").Value = -3162.536599999999
$ws.Range("H135").Value = 15626782
$ws.Range("I135").Value = 16130840
$ws.Range("K135").Value = 145177560
$ws.Range("M135").Value = -145175025
$ws.Range("H137").Value = 2900.6775
$ws.Range("I137").Value = 2708.8845
$ws.Range("J137").Value = 3898
$ws.Range("K137").Value = 8126.6535
$ws.Range("L137").Value = 11694
$ws.Range("M137").Value = -5576.6535
$ws.Range("N137").Value = -16794
$ws.Range("H138").Value = 3822.8108
$ws.Range("I138").Value = 2110.7778
$ws.Range("J138").Value = 4373.107
$ws.Range("K138").Value = 6332.3334
$ws.Range("L138").Value = 13119.321
$ws.Range("M138").Value = -1192.3334
$ws.Range("N138").Value = -23399.321
$ws.Range("H141").Value = 3948.8
$ws.Range("I141").Value = 2863.3333
$ws.Range("K141").Value = 8589.999899999999
$ws.Range("M141").Value = -3409.999899999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 20333.705
$ws.Range("I2").Value = 22778.2
$ws.Range("K2").Value = 22778.2
$ws.Range("M2").Value = -22665.2
$ws.Range("H5").Value = 500.5
$ws.Range("I5").Value = 500.5
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 500.5
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -388.5
$ws.Range("N5").ClearContents()
$ws.Range("H45").Value = 9012.799999999999
$ws.Range("I45").Value = 8675
$ws.Range("J45").Value = 9238
$ws.Range("K45").Value = 8675
$ws.Range("L45").Value = 9238
$ws.Range("M45").Value = -8298
$ws.Range("N45").Value = -9992
$ws.Range("H63").Value = 9618.799999999999
$ws.Range("I63").Value = 5169.7144
$ws.Range("K63").Value = 5169.7144
$ws.Range("M63").Value = -4483.7144
$ws.Range("H66").Value = 9618.799999999999
$ws.Range("I66").Value = 5169.7144
$ws.Range("K66").Value = 25848.572
$ws.Range("M66").Value = -22416.572
$ws.Range("H88").Value = 1379.6666
$ws.Range("J88").Value = 1374.5834
$ws.Range("L88").Value = 1374.5834
$ws.Range("N88").Value = -2186.5834
$ws.Range("H91").Value = 1379.6666
$ws.Range("J91").Value = 1374.5834
$ws.Range("L91").Value = 1374.5834
$ws.Range("N91").Value = -4182.5834
$ws.Range("H97").Value = 1175.9459
$ws.Range("I97").Value = 642.82855
$ws.Range("K97").Value = 642.82855
$ws.Range("M97").Value = -146.82855
$ws.Range("H116").Value = 20333.705
$ws.Range("I116").Value = 22778.2
$ws.Range("K116").Value = 22778.2
$ws.Range("M116").Value = -20484.2
$ws.Range("H132").Value = 2726.9092
$ws.Range("I132").Value = 2778.25
$ws.Range("K132").Value = 8334.75
$ws.Range("M132").Value = -5804.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 20333.705
$ws.Range("I3").Value = 22778.2
$ws.Range("K3").Value = 22778.2
$ws.Range("M3").Value = -22664.2
$ws.Range("H4").Value = 500.5
$ws.Range("I4").Value = 500.5
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 500.5
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -385.5
$ws.Range("N4").ClearContents()
$ws.Range("H134").Value = 1677.1063
$ws.Range("I134").Value = 1301.8049
$ws.Range("J134").Value = 4241.6665
$ws.Range("K134").Value = 3905.4147
$ws.Range("L134").Value = 12724.9995
$ws.Range("M134").Value = -1370.4147
$ws.Range("N134").Value = -17794.9995
$ws.Range("H140").Value = 70718.47

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 129.22223
$ws.Range("I7").Value = 75
$ws.Range("J7").Value = 156.33333
$ws.Range("K7").Value = 75
$ws.Range("L7").Value = 156.33333
$ws.Range("M7").Value = 38
$ws.Range("N7").Value = -382.33333
$ws.Range("H13").Value = 69.75
$ws.Range("J13").Value = 69.75
$ws.Range("L13").Value = 69.75
$ws.Range("N13").Value = -347.75
$ws.Range("H21").Value = 1059
$ws.Range("I21").Value = 1059
$ws.Range("K21").Value = 1059
$ws.Range("M21").Value = -824
$ws.Range("H31").Value = 2280.724
$ws.Range("I31").Value = 2309.889
$ws.Range("J31").Value = 2233
$ws.Range("K31").Value = 2309.889
$ws.Range("L31").Value = 2233
$ws.Range("M31").Value = -2014.889
$ws.Range("N31").Value = -2823
$ws.Range("H34").Value = 2280.724
$ws.Range("I34").Value = 2309.889
$ws.Range("J34").Value = 2233
$ws.Range("K34").Value = 2309.889
$ws.Range("L34").Value = 2233
$ws.Range("M34").Value = -2107.889
$ws.Range("N34").Value = -2637
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H58").Value = 2582.375
$ws.Range("I58").Value = 1562.1875
$ws.Range("K58").Value = 1562.1875
$ws.Range("M58").Value = -1359.1875
$ws.Range("H99").Value = 22446148
$ws.Range("J99").Value = 22230136
$ws.Range("L99").Value = 22230136
$ws.Range("N99").Value = -22233132
$ws.Range("H123").Value = 49999.09
$ws.Range("J123").Value = 49999.09
$ws.Range("L123").Value = 49999.09
$ws.Range("N123").Value = -59799.09
$ws.Range("H126").Value = 22446148
$ws.Range("J126").Value = 22230136
$ws.Range("L126").Value = 66690408
$ws.Range("N126").Value = -66695348
$ws.Range("H132").Value = 2226.7368
$ws.Range("I132").Value = 2194.5881
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 6583.7643
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -4053.7643
$ws.Range("N132").Value = -12560
$ws.Range("H134").Value = 2891.4375
$ws.Range("I134").Value = 2351.64
$ws.Range("K134").Value = 7054.92
$ws.Range("M134").Value = -4519.92
$ws.Range("H136").Value = 2582.375
$ws.Range("I136").Value = 1562.1875
$ws.Range("K136").Value = 4686.5625
$ws.Range("M136").Value = -2136.5625

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 127
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 127
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 381
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -719
$ws.Range("H36").Value = 1220.5454
$ws.Range("I36").Value = 1220.5454
$ws.Range("K36").Value = 3661.6362
$ws.Range("M36").Value = -3492.6362
$ws.Range("H44").Value = 437.5
$ws.Range("I44").Value = 466.66666
$ws.Range("K44").Value = 1399.99998
$ws.Range("M44").Value = -1001.99998
$ws.Range("H141").Value = 142864690
$ws.Range("I141").Value = 142864690
$ws.Range("K141").Value = 428594070
$ws.Range("M141").Value = -428588890

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2968.2
$ws.Range("I126").Value = 2968.2
$ws.Range("K126").Value = 8904.599999999999
$ws.Range("M126").Value = -6434.599999999999
$ws.Range("H132").Value = 2633.5833
$ws.Range("I132").Value = 2339.875
$ws.Range("J132").Value = 3221
$ws.Range("K132").Value = 7019.625
$ws.Range("L132").Value = 9663
$ws.Range("M132").Value = -4489.625
$ws.Range("N132").Value = -14723

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 12025.091
$ws.Range("I40").Value = 13569
$ws.Range("K40").Value = 13569
$ws.Range("M40").Value = -13433
$ws.Range("H61").Value = 2706.4614
$ws.Range("I61").Value = 2120.7
$ws.Range("K61").Value = 2120.7
$ws.Range("M61").Value = -1918.7
$ws.Range("H93").Value = 6700
$ws.Range("I93").Value = 7933.3335
$ws.Range("J93").Value = 3000
$ws.Range("K93").Value = 7933.3335
$ws.Range("L93").Value = 3000
$ws.Range("M93").Value = -6685.3335
$ws.Range("N93").Value = -5496
$ws.Range("H113").Value = 2706.4614
$ws.Range("I113").Value = 2120.7
$ws.Range("K113").Value = 2120.7
$ws.Range("M113").Value = 49.30000000000018

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3480.8
$ws.Range("I126").Value = 2702
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 8106
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -5636
$ws.Range("N126").Value = -16940
